$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set column D width similar to A:B (12.77734375 characters)
$ws.Columns.Item(4).ColumnWidth = 12

# Add new column D values (shared string lookups added by the edit)
# Order chosen so the shared-strings table is built in the same sequence
# as found in the target workbook (ZC05, SI001, ZW06, ZC04, ZAGNG0301B).
$ws.Range("D12").Value = "ZC05"
$ws.Range("D15").Value = "SI001"
$ws.Range("D1").Value = "ZW06"
$ws.Range("D14").Value = "ZC04"
$ws.Range("D6").Value = "ZAGNG0301B"
$ws.Range("D13").Value = "ZC05"

# Page setup change (portrait, paper size 9 = A4)
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
